$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 570321
$ws.Range("D2").Value = 153186
$ws.Range("E2").Value = 948706158

$ws.Range("C8").Value = 2482
$ws.Range("E8").Value = 13496822

$ws.Range("C10").Value = 236317
$ws.Range("D10").Value = 62596
$ws.Range("E10").Value = 925634224

$ws.Range("C13").Value = 125242
$ws.Range("E13").Value = 515139004

$ws.Range("C16").Value = 6868
$ws.Range("E16").Value = 14631302

$ws.Range("C19").Value = 16322
$ws.Range("E19").Value = 51969578

$ws.Range("C21").Value = 133596
$ws.Range("E21").Value = 222045099

$ws.Range("C26").Value = 773
$ws.Range("E26").Value = 3357284

$ws.Range("C27").Value = 62874
$ws.Range("E27").Value = 232942074

$ws.Range("C30").Value = 23830
$ws.Range("E30").Value = 92814778

$ws.Range("C35").Value = 4172
$ws.Range("E35").Value = 12717362

$ws.Range("C36").Value = 161981
$ws.Range("E36").Value = 280371594

$ws.Range("C39").Value = 2933
$ws.Range("E39").Value = 16306741

$ws.Range("C41").Value = 88623
$ws.Range("D41").Value = 23788
$ws.Range("E41").Value = 351046253

$ws.Range("C44").Value = 19623
$ws.Range("E44").Value = 89338848

$ws.Range("C47").Value = 5395
$ws.Range("E47").Value = 17851761

$ws.Range("C48").Value = 116547
$ws.Range("E48").Value = 195641997

$ws.Range("C54").Value = 53433
$ws.Range("D54").Value = 14602
$ws.Range("E54").Value = 192399888

$ws.Range("C57").Value = 21935
$ws.Range("D57").Value = 5809
$ws.Range("E57").Value = 79366993

$ws.Range("C58").Value = 2349
$ws.Range("E58").Value = 6216424

$ws.Range("C60").Value = 3566
$ws.Range("E60").Value = 10334591

$ws.Range("C62").Value = 36234
$ws.Range("E62").Value = 67095788

$ws.Range("C66").Value = 17562
$ws.Range("E66").Value = 81306547

$ws.Range("C68").Value = 11532
$ws.Range("E68").Value = 47936777

$ws.Range("C71").Value = 248777
$ws.Range("E71").Value = 431778861

$ws.Range("C75").Value = 1281
$ws.Range("E75").Value = 6314020

$ws.Range("C77").Value = 123801
$ws.Range("D77").Value = 33947
$ws.Range("E77").Value = 459954564

$ws.Range("C80").Value = 59582
$ws.Range("D80").Value = 16035
$ws.Range("E80").Value = 230144387

$ws.Range("C82").Value = 237
$ws.Range("E82").Value = 2300609

$ws.Range("C83").Value = 13792
$ws.Range("D83").Value = 6122
$ws.Range("E83").Value = 65663750

$ws.Range("C85").Value = 6629
$ws.Range("E85").Value = 20650262

$ws.Range("C86").Value = 50277
$ws.Range("D86").Value = 12014
$ws.Range("E86").Value = 77769344

$ws.Range("C89").Value = 11918
$ws.Range("E89").Value = 24930546

$ws.Range("C91").Value = 11137
$ws.Range("E91").Value = 21303261

$ws.Range("C94").Value = 20851
$ws.Range("E94").Value = 42614422

$ws.Range("C95").Value = 4597
$ws.Range("E95").Value = 11019559

$ws.Range("C97").Value = 6688
$ws.Range("E97").Value = 15706386

$ws.Range("C101").Value = 242712
$ws.Range("D101").Value = 66889
$ws.Range("E101").Value = 394924462

$ws.Range("C106").Value = 2747
$ws.Range("E106").Value = 13350208

$ws.Range("C108").Value = 100472
$ws.Range("E108").Value = 364851466

$ws.Range("C111").Value = 54910
$ws.Range("D111").Value = 13883
$ws.Range("E111").Value = 201934250

$ws.Range("C113").Value = 2593
$ws.Range("D113").Value = 1015
$ws.Range("E113").Value = 6859754

$ws.Range("C114").Value = 5260
$ws.Range("E114").Value = 15302843

$ws.Range("C116").Value = 975953
$ws.Range("D116").Value = 216387
$ws.Range("E116").Value = 1658027331

$ws.Range("C117").Value = 201
$ws.Range("E117").Value = 497029

$ws.Range("C118").Value = 486
$ws.Range("E118").Value = 2009598

$ws.Range("C121").Value = 4832
$ws.Range("E121").Value = 40297131

$ws.Range("C123").Value = 420787
$ws.Range("D123").Value = 101956
$ws.Range("E123").Value = 1563900630

$ws.Range("C124").Value = 1920
$ws.Range("E124").Value = 31400967

$ws.Range("C126").Value = 388862
$ws.Range("D126").Value = 86777
$ws.Range("E126").Value = 1403952855

$ws.Range("C128").Value = 4909
$ws.Range("E128").Value = 9350400

$ws.Range("C130").Value = 15779
$ws.Range("E130").Value = 48942710

$ws.Range("C133").Value = 60584
$ws.Range("E133").Value = 87334381

$ws.Range("C138").Value = 17861
$ws.Range("E138").Value = 36143156

$ws.Range("C140").Value = 5068
$ws.Range("E140").Value = 10032581

$ws.Range("C145").Value = 27871
$ws.Range("E145").Value = 42864758

$ws.Range("C148").Value = 11554
$ws.Range("D148").Value = 3133
$ws.Range("E148").Value = 28395346

$ws.Range("C150").Value = 8221
$ws.Range("D150").Value = 2041
$ws.Range("E150").Value = 18018767

$ws.Range("C153").Value = 37846
$ws.Range("E153").Value = 92194203

$ws.Range("C154").Value = 3593
$ws.Range("D154").Value = 665
$ws.Range("E154").Value = 9301360

$ws.Range("C158").Value = 149370
$ws.Range("E158").Value = 256715678

$ws.Range("C163").Value = 2093
$ws.Range("E163").Value = 11762191

$ws.Range("C165").Value = 66673
$ws.Range("E165").Value = 256580669

$ws.Range("C167").Value = 27109
$ws.Range("E167").Value = 109979837

$ws.Range("C169").Value = 1978
$ws.Range("E169").Value = 4170328

$ws.Range("C170").Value = 4101
$ws.Range("D170").Value = 1076
$ws.Range("E170").Value = 12589862

$ws.Range("C171").Value = 396408
$ws.Range("D171").Value = 113384
$ws.Range("E171").Value = 632151430

$ws.Range("C179").Value = 165565
$ws.Range("E179").Value = 618732270

$ws.Range("C182").Value = 67436
$ws.Range("E182").Value = 262137774

$ws.Range("C185").Value = 9378
$ws.Range("D185").Value = 3560
$ws.Range("E185").Value = 28314614

$ws.Range("C187").Value = 11213
$ws.Range("E187").Value = 31271976

$ws.Range("C189").Value = 454754
$ws.Range("D189").Value = 124900
$ws.Range("E189").Value = 702772305

$ws.Range("C197").Value = 186050
$ws.Range("D197").Value = 48724
$ws.Range("E197").Value = 674921907

$ws.Range("C200").Value = 105369
$ws.Range("D200").Value = 26188
$ws.Range("E200").Value = 376522395

$ws.Range("C203").Value = 7935
$ws.Range("E203").Value = 17962646

$ws.Range("C206").Value = 14147
$ws.Range("E206").Value = 38375716

$ws.Range("C208").Value = 177356
$ws.Range("E208").Value = 292210037

$ws.Range("C214").Value = 94969
$ws.Range("D214").Value = 26659
$ws.Range("E214").Value = 356961634

$ws.Range("C217").Value = 26206
$ws.Range("E217").Value = 110294486

$ws.Range("C221").Value = 5703
$ws.Range("E221").Value = 17623971

$ws.Range("C222").Value = 463130
$ws.Range("D222").Value = 119162
$ws.Range("E222").Value = 740159957

$ws.Range("C228").Value = 2607
$ws.Range("E228").Value = 15663545

$ws.Range("C230").Value = 196294
$ws.Range("D230").Value = 48719
$ws.Range("E230").Value = 747943129

$ws.Range("C231").Value = 371
$ws.Range("E231").Value = 7014089

$ws.Range("C233").Value = 140500
$ws.Range("D233").Value = 32886
$ws.Range("E233").Value = 518632819

$ws.Range("C236").Value = 4636
$ws.Range("D236").Value = 1593
$ws.Range("E236").Value = 11765069

$ws.Range("C239").Value = 11463
$ws.Range("E239").Value = 32799926
